$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$nbsp = [char]0x00A0

# Ensure the full target range is formatted as Text so numeric-looking strings
# (e.g. "13", "0.00") are preserved as text instead of being auto-converted to numbers.
$ws.Range("A1:K12").NumberFormat = "@"

# Row 1
$ws.Cells.Item(1, 1).Value2 = "venue"
$ws.Cells.Item(1, 2).Value2 = "date"
$ws.Cells.Item(1, 3).Value2 = "result"
$ws.Cells.Item(1, 4).Value2 = "ownTeam"
$ws.Cells.Item(1, 5).Value2 = "oppTeam"
$ws.Cells.Item(1, 6).Value2 = "batsman"
$ws.Cells.Item(1, 7).Value2 = "totalRuns"
$ws.Cells.Item(1, 8).Value2 = "totalBalls"
$ws.Cells.Item(1, 9).Value2 = "total4s"
$ws.Cells.Item(1, 10).Value2 = "total6s"
$ws.Cells.Item(1, 11).Value2 = "sr"

# Row 2
$ws.Cells.Item(2, 1).Value2 = " Sharjah"
$ws.Cells.Item(2, 2).Value2 = " September 27 2020"
$ws.Cells.Item(2, 3).Value2 = "Royals won by 4 wickets (with 3 balls remaining)"
$ws.Cells.Item(2, 4).Value2 = "Kings XI Punjab"
$ws.Cells.Item(2, 5).Value2 = "Rajasthan Royals"
$ws.Cells.Item(2, 6).Value2 = "Glenn Maxwell" + $nbsp
$ws.Cells.Item(2, 7).Value2 = "13"
$ws.Cells.Item(2, 8).Value2 = "9"
$ws.Cells.Item(2, 9).Value2 = "2"
$ws.Cells.Item(2, 10).Value2 = "0"
$ws.Cells.Item(2, 11).Value2 = "144.44"

# Row 3
$ws.Cells.Item(3, 1).Value2 = " Abu Dhabi"
$ws.Cells.Item(3, 2).Value2 = " October 30 2020"
$ws.Cells.Item(3, 3).Value2 = "Royals won by 7 wickets (with 15 balls remaining)"
$ws.Cells.Item(3, 4).Value2 = "Kings XI Punjab"
$ws.Cells.Item(3, 5).Value2 = "Rajasthan Royals"
$ws.Cells.Item(3, 6).Value2 = "Glenn Maxwell" + $nbsp
$ws.Cells.Item(3, 7).Value2 = "6"
$ws.Cells.Item(3, 8).Value2 = "6"
$ws.Cells.Item(3, 9).Value2 = "1"
$ws.Cells.Item(3, 10).Value2 = "0"
$ws.Cells.Item(3, 11).Value2 = "100.00"

# Row 4
$ws.Cells.Item(4, 1).Value2 = " Dubai (DSC)"
$ws.Cells.Item(4, 2).Value2 = " October 04 2020"
$ws.Cells.Item(4, 3).Value2 = "Super Kings won by 10 wickets (with 14 balls remaining)"
$ws.Cells.Item(4, 4).Value2 = "Kings XI Punjab"
$ws.Cells.Item(4, 5).Value2 = "Chennai Super Kings"
$ws.Cells.Item(4, 6).Value2 = "Glenn Maxwell" + $nbsp
$ws.Cells.Item(4, 7).Value2 = "11"
$ws.Cells.Item(4, 8).Value2 = "7"
$ws.Cells.Item(4, 9).Value2 = "1"
$ws.Cells.Item(4, 10).Value2 = "0"
$ws.Cells.Item(4, 11).Value2 = "157.14"

# Row 5
$ws.Cells.Item(5, 1).Value2 = " Dubai (DSC)"
$ws.Cells.Item(5, 2).Value2 = " October 24 2020"
$ws.Cells.Item(5, 3).Value2 = "Kings XI won by 12 runs"
$ws.Cells.Item(5, 4).Value2 = "Kings XI Punjab"
$ws.Cells.Item(5, 5).Value2 = "Sunrisers Hyderabad"
$ws.Cells.Item(5, 6).Value2 = "Glenn Maxwell" + $nbsp
$ws.Cells.Item(5, 7).Value2 = "12"
$ws.Cells.Item(5, 8).Value2 = "13"
$ws.Cells.Item(5, 9).Value2 = "0"
$ws.Cells.Item(5, 10).Value2 = "0"
$ws.Cells.Item(5, 11).Value2 = "92.30"

# Row 6
$ws.Cells.Item(6, 1).Value2 = " Dubai (DSC)"
$ws.Cells.Item(6, 2).Value2 = " September 24 2020"
$ws.Cells.Item(6, 3).Value2 = "Kings XI won by 97 runs"
$ws.Cells.Item(6, 4).Value2 = "Kings XI Punjab"
$ws.Cells.Item(6, 5).Value2 = "Royal Challengers Bangalore"
$ws.Cells.Item(6, 6).Value2 = "Glenn Maxwell" + $nbsp
$ws.Cells.Item(6, 7).Value2 = "5"
$ws.Cells.Item(6, 8).Value2 = "6"
$ws.Cells.Item(6, 9).Value2 = "0"
$ws.Cells.Item(6, 10).Value2 = "0"
$ws.Cells.Item(6, 11).Value2 = "83.33"

# Row 7
$ws.Cells.Item(7, 1).Value2 = " Abu Dhabi"
$ws.Cells.Item(7, 2).Value2 = " October 01 2020"
$ws.Cells.Item(7, 3).Value2 = "Mumbai won by 48 runs"
$ws.Cells.Item(7, 4).Value2 = "Kings XI Punjab"
$ws.Cells.Item(7, 5).Value2 = "Mumbai Indians"
$ws.Cells.Item(7, 6).Value2 = "Glenn Maxwell" + $nbsp
$ws.Cells.Item(7, 7).Value2 = "11"
$ws.Cells.Item(7, 8).Value2 = "18"
$ws.Cells.Item(7, 9).Value2 = "0"
$ws.Cells.Item(7, 10).Value2 = "0"
$ws.Cells.Item(7, 11).Value2 = "61.11"

# Row 8
$ws.Cells.Item(8, 1).Value2 = " Abu Dhabi"
$ws.Cells.Item(8, 2).Value2 = " October 10 2020"
$ws.Cells.Item(8, 3).Value2 = "KKR won by 2 runs"
$ws.Cells.Item(8, 4).Value2 = "Kings XI Punjab"
$ws.Cells.Item(8, 5).Value2 = "Kolkata Knight Riders"
$ws.Cells.Item(8, 6).Value2 = "Glenn Maxwell" + $nbsp
$ws.Cells.Item(8, 7).Value2 = "10"
$ws.Cells.Item(8, 8).Value2 = "5"
$ws.Cells.Item(8, 9).Value2 = "2"
$ws.Cells.Item(8, 10).Value2 = "0"
$ws.Cells.Item(8, 11).Value2 = "200.00"

# Row 9
$ws.Cells.Item(9, 1).Value2 = " Dubai (DSC)"
$ws.Cells.Item(9, 2).Value2 = " October 20 2020"
$ws.Cells.Item(9, 3).Value2 = "Kings XI won by 5 wickets (with 6 balls remaining)"
$ws.Cells.Item(9, 4).Value2 = "Kings XI Punjab"
$ws.Cells.Item(9, 5).Value2 = "Delhi Capitals"
$ws.Cells.Item(9, 6).Value2 = "Glenn Maxwell" + $nbsp
$ws.Cells.Item(9, 7).Value2 = "32"
$ws.Cells.Item(9, 8).Value2 = "24"
$ws.Cells.Item(9, 9).Value2 = "3"
$ws.Cells.Item(9, 10).Value2 = "0"
$ws.Cells.Item(9, 11).Value2 = "133.33"

# Row 10
$ws.Cells.Item(10, 1).Value2 = " Dubai (DSC)"
$ws.Cells.Item(10, 2).Value2 = " September 20 2020"
$ws.Cells.Item(10, 3).Value2 = "Match tied (Capitals won the one-over eliminator)"
$ws.Cells.Item(10, 4).Value2 = "Kings XI Punjab"
$ws.Cells.Item(10, 5).Value2 = "Delhi Capitals"
$ws.Cells.Item(10, 6).Value2 = "Glenn Maxwell" + $nbsp
$ws.Cells.Item(10, 7).Value2 = "1"
$ws.Cells.Item(10, 8).Value2 = "4"
$ws.Cells.Item(10, 9).Value2 = "0"
$ws.Cells.Item(10, 10).Value2 = "0"
$ws.Cells.Item(10, 11).Value2 = "25.00"

# Row 11
$ws.Cells.Item(11, 1).Value2 = " Dubai (DSC)"
$ws.Cells.Item(11, 2).Value2 = " October 18 2020"
$ws.Cells.Item(11, 3).Value2 = "Match tied (Kings XI won the one-over eliminator)"
$ws.Cells.Item(11, 4).Value2 = "Kings XI Punjab"
$ws.Cells.Item(11, 5).Value2 = "Mumbai Indians"
$ws.Cells.Item(11, 6).Value2 = "Glenn Maxwell" + $nbsp
$ws.Cells.Item(11, 7).Value2 = "0"
$ws.Cells.Item(11, 8).Value2 = "2"
$ws.Cells.Item(11, 9).Value2 = "0"
$ws.Cells.Item(11, 10).Value2 = "0"
$ws.Cells.Item(11, 11).Value2 = "0.00"

# Row 12
$ws.Cells.Item(12, 1).Value2 = " Dubai (DSC)"
$ws.Cells.Item(12, 2).Value2 = " October 08 2020"
$ws.Cells.Item(12, 3).Value2 = "Sunrisers won by 69 runs"
$ws.Cells.Item(12, 4).Value2 = "Kings XI Punjab"
$ws.Cells.Item(12, 5).Value2 = "Sunrisers Hyderabad"
$ws.Cells.Item(12, 6).Value2 = "Glenn Maxwell" + $nbsp
$ws.Cells.Item(12, 7).Value2 = "7"
$ws.Cells.Item(12, 8).Value2 = "12"
$ws.Cells.Item(12, 9).Value2 = "0"
$ws.Cells.Item(12, 10).Value2 = "0"
$ws.Cells.Item(12, 11).Value2 = "58.33"
